$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.290.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.678.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '674.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.16'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.45%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.146'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.91'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.435'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000231'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.293.36'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.35'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.673.65'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.229.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.40'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.73'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.821.81'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000122'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.88'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.08'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.67'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.75'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.60'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.32%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.96'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.98'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.669.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.49%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.19'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.22'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.12%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0902'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '170.05'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.941'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.64'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.44%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000276'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.67'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.08'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.78'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.05%  '
